$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 332
$ws.Range("I41").Value = 346.25
$ws.Range("J41").Value = 275
$ws.Range("K41").Value = 346.25
$ws.Range("L41").Value = 275
$ws.Range("M41").Value = 93.75
$ws.Range("N41").Value = -1155
$ws.Range("H43").Value = 3965.6667
$ws.Range("I43").Value = 2099
$ws.Range("K43").Value = 2099
$ws.Range("M43").Value = -2030
$ws.Range("H51").Value = 4233
$ws.Range("I51").Value = 4233
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 4233
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -3749
$ws.Range("H53").Value = 84241.836
$ws.Range("I53").Value = 1512.3334
$ws.Range("J53").Value = 166971.33
$ws.Range("K53").Value = 1512.3334
$ws.Range("L53").Value = 166971.33
$ws.Range("M53").Value = -875.3334
$ws.Range("N53").Value = -168245.33
$ws.Range("H61").Value = 25
$ws.Range("I61").Value = 25
$ws.Range("K61").Value = 75
$ws.Range("M61").Value = 97
$ws.Range("H64").Value = 9169.286
$ws.Range("J64").Value = 11888.889
$ws.Range("L64").Value = 11888.889
$ws.Range("N64").Value = -12384.889
$ws.Range("H67").Value = 9169.286
$ws.Range("J67").Value = 11888.889
$ws.Range("L67").Value = 11888.889
$ws.Range("N67").Value = -13604.889
$ws.Range("H86").Value = 2899
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2899
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H113").Value = 166670370
$ws.Range("J113").Value = 5006
$ws.Range("L113").Value = 5006
$ws.Range("N113").Value = -11514
$ws.Range("H135").Value = 1575.7858
$ws.Range("J135").Value = 1497.5
$ws.Range("L135").Value = 13477.5
$ws.Range("N135").Value = -18547.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 73.125
$ws.Range("I5").Value = 65.71429000000001
$ws.Range("K5").Value = 65.71429000000001
$ws.Range("M5").Value = 46.28570999999999
$ws.Range("H32").Value = 2137
$ws.Range("I32").Value = 2003.7288
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 2003.7288
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -1716.7288
$ws.Range("N32").Value = -10574
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").ClearContents()
$ws.Range("N37").Value = 0
$ws.Range("H74").Value = 3837.12
$ws.Range("I74").Value = 2866.4348
$ws.Range("K74").Value = 2866.4348
$ws.Range("M74").Value = -1992.4348
$ws.Range("H77").Value = 3837.12
$ws.Range("I77").Value = 2866.4348
$ws.Range("K77").Value = 14332.174
$ws.Range("M77").Value = -9964.173999999999
$ws.Range("H97").Value = 820.6667
$ws.Range("I97").Value = 647.6667
$ws.Range("K97").Value = 647.6667
$ws.Range("M97").Value = -151.6667
$ws.Range("H102").Value = 5118.375
$ws.Range("J102").Value = 8250
$ws.Range("L102").Value = 8250
$ws.Range("N102").Value = -11494

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 73.125
$ws.Range("I4").Value = 65.71429000000001
$ws.Range("K4").Value = 65.71429000000001
$ws.Range("M4").Value = 49.28570999999999
$ws.Range("H46").Value = 30000
$ws.Range("J46").Value = 30000
$ws.Range("L46").Value = 30000
$ws.Range("N46").Value = -30596
$ws.Range("H99").Value = 3918.641
$ws.Range("I99").Value = 2699.077
$ws.Range("J99").Value = 6357.769
$ws.Range("K99").Value = 2699.077
$ws.Range("L99").Value = 6357.769
$ws.Range("M99").Value = -1201.077
$ws.Range("N99").Value = -9353.769
$ws.Range("H134").Value = 3785.9487
$ws.Range("I134").Value = 3785.9487
$ws.Range("K134").Value = 11357.8461
$ws.Range("M134").Value = -8822.846099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4142.7144
$ws.Range("I16").Value = 2931
$ws.Range("J16").Value = 5354.4287
$ws.Range("K16").Value = 2931
$ws.Range("L16").Value = 5354.4287
$ws.Range("M16").Value = -2644
$ws.Range("N16").Value = -5928.4287
$ws.Range("H31").Value = 6006
$ws.Range("I31").Value = 4571.2666
$ws.Range("K31").Value = 4571.2666
$ws.Range("M31").Value = -4276.2666
$ws.Range("H34").Value = 6006
$ws.Range("I34").Value = 4571.2666
$ws.Range("K34").Value = 4571.2666
$ws.Range("M34").Value = -4369.2666
$ws.Range("H94").Value = 4829.5713
$ws.Range("J94").Value = 4829.5713
$ws.Range("L94").Value = 4829.5713
$ws.Range("N94").Value = -5731.5713
$ws.Range("H113").Value = 4142.7144
$ws.Range("I113").Value = 2931
$ws.Range("J113").Value = 5354.4287
$ws.Range("K113").Value = 2931
$ws.Range("L113").Value = 5354.4287
$ws.Range("M113").Value = -761
$ws.Range("N113").Value = -9694.4287
$ws.Range("H132").Value = 2244.111
$ws.Range("I132").Value = 1899.875
$ws.Range("K132").Value = 5699.625
$ws.Range("M132").Value = -3169.625
$ws.Range("H135").Value = 93467.60000000001
$ws.Range("J135").Value = 93467.60000000001
$ws.Range("L135").Value = 93467.60000000001
$ws.Range("N135").Value = -103607.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 463.875
$ws.Range("I10").Value = 49.304348
$ws.Range("K10").Value = 147.913044
$ws.Range("M10").Value = -8.913043999999985
$ws.Range("H18").Value = 2257.25
$ws.Range("I18").Value = 1515
$ws.Range("K18").Value = 4545
$ws.Range("M18").Value = -4376
$ws.Range("H121").Value = 20000598
$ws.Range("I121").Value = 666.3333
$ws.Range("J121").Value = 50000496
$ws.Range("K121").Value = 1998.9999
$ws.Range("L121").Value = 150001488
$ws.Range("M121").Value = -688.9999
$ws.Range("N121").Value = -150004108
$ws.Range("H132").Value = 41667824
$ws.Range("I132").Value = 125000470
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 1125004230
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -1125001700
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 17519.8
$ws.Range("I22").Value = 9199.666999999999
$ws.Range("K22").Value = 9199.666999999999
$ws.Range("M22").Value = -8670.666999999999
$ws.Range("H102").Value = 2800.3215
$ws.Range("J102").Value = 4368.8
$ws.Range("L102").Value = 4368.8
$ws.Range("N102").Value = -7612.8
$ws.Range("H104").Value = 39995
$ws.Range("J104").Value = 39995
$ws.Range("L104").Value = 39995
$ws.Range("N104").Value = -46983
$ws.Range("H122").Value = 2439
$ws.Range("I122").Value = 2439
$ws.Range("K122").Value = 7317
$ws.Range("M122").Value = -4867
$ws.Range("H126").Value = 4802.222
$ws.Range("I126").Value = 5642.7144
$ws.Range("K126").Value = 16928.1432
$ws.Range("M126").Value = -14458.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2567.25
$ws.Range("I16").Value = 2843.6
$ws.Range("K16").Value = 2843.6
$ws.Range("M16").Value = -2673.6
$ws.Range("H22").Value = 2969
$ws.Range("J22").Value = 4333.3335
$ws.Range("L22").Value = 4333.3335
$ws.Range("N22").Value = -4923.3335
$ws.Range("H27").Value = 2969
$ws.Range("J27").Value = 4333.3335
$ws.Range("L27").Value = 4333.3335
$ws.Range("N27").Value = -4547.3335
$ws.Range("H40").Value = 2997.4119
$ws.Range("I40").Value = 2996.75
$ws.Range("K40").Value = 2996.75
$ws.Range("M40").Value = -2860.75
$ws.Range("H55").Value = 490.4138
$ws.Range("I55").Value = 592.7368
$ws.Range("K55").Value = 592.7368
$ws.Range("M55").Value = -419.7368
$ws.Range("H100").Value = 2781525.8
$ws.Range("J100").Value = 4999.5713
$ws.Range("L100").Value = 4999.5713
$ws.Range("N100").Value = -6081.5713
$ws.Range("H132").Value = 12606.667
$ws.Range("I132").Value = 15125.125
$ws.Range("J132").Value = 7569.75
$ws.Range("K132").Value = 45375.375
$ws.Range("L132").Value = 22709.25
$ws.Range("M132").Value = -42845.375
$ws.Range("N132").Value = -27769.25
$ws.Range("H136").Value = 4916.8
$ws.Range("J136").Value = 4833.25
$ws.Range("L136").Value = 14499.75
$ws.Range("N136").Value = -19599.75
$ws.Range("H139").Value = 83928.336
$ws.Range("J139").Value = 83928.336
$ws.Range("L139").Value = 83928.336
$ws.Range("N139").Value = -94208.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5439.885
$ws.Range("I126").Value = 5682.087
$ws.Range("J126").Value = 3583
$ws.Range("K126").Value = 17046.261
$ws.Range("L126").Value = 10749
$ws.Range("M126").Value = -14576.261
$ws.Range("N126").Value = -15689
$ws.Range("H132").Value = 5947.7427
$ws.Range("I132").Value = 5306.615
$ws.Range("J132").Value = 7799.8887
$ws.Range("K132").Value = 15919.845
$ws.Range("L132").Value = 23399.6661
$ws.Range("M132").Value = -13389.845
$ws.Range("N132").Value = -28459.6661
$ws.Range("H136").Value = 5320.9473
$ws.Range("I136").Value = 4422.552
$ws.Range("K136").Value = 13267.656
$ws.Range("M136").Value = -10717.656
